$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains exact text formatting (avoid Excel auto-converting
# numeric-looking strings like "10.40" or "1.71" into numbers and dropping
# trailing zeros / formatting).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '60.123.96'
$ws.Range("E2").Value = '  +2.63%  '
$ws.Range("D3").Value = '2.630.15'
$ws.Range("E3").Value = '  +0.16%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '566.91'
$ws.Range("E5").Value = '  +6.01%  '
$ws.Range("D6").Value = '145.43'
$ws.Range("E6").Value = '  +2.11%  '
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.24%  '
$ws.Range("D8").Value = '0.608'
$ws.Range("E8").Value = '  +7.13%  '
$ws.Range("E9").Value = '  -2.92%  '
$ws.Range("D10").Value = '0.104'
$ws.Range("E10").Value = '  +4.03%  '
$ws.Range("E11").Value = '  +6.35%  '
$ws.Range("E12").Value = '  +2.51%  '
$ws.Range("D13").Value = '3.096.45'
$ws.Range("E13").Value = '  +0.30%  '
$ws.Range("D14").Value = '60.115.04'
$ws.Range("E14").Value = '  +2.72%  '
$ws.Range("D15").Value = '21.65'
$ws.Range("E15").Value = '  +3.76%  '
$ws.Range("D16").Value = '2.647.55'
$ws.Range("E16").Value = '  +0.69%  '
$ws.Range("D17").Value = '0.0000136'
$ws.Range("E17").Value = '  +2.59%  '
$ws.Range("D18").Value = '4.59'
$ws.Range("E18").Value = '  +4.63%  '
$ws.Range("D19").Value = '343.27'
$ws.Range("E19").Value = '  +2.46%  '
$ws.Range("D20").Value = '10.40'
$ws.Range("E20").Value = '  +2.37%  '
$ws.Range("D21").Value = '6.27'
$ws.Range("E21").Value = '  +0.95%  '
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.13%  '
$ws.Range("D23").Value = '66.54'
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("E24").Value = '  +5.04%  '
$ws.Range("D25").Value = '0.166'
$ws.Range("E25").Value = '  +1.79%  '
$ws.Range("D26").Value = '0.996'
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("D27").Value = '7.29'
$ws.Range("E27").Value = '  +2.44%  '
$ws.Range("D28").Value = '0.0₃0772'
$ws.Range("E28").Value = '  +4.69%  '
$ws.Range("D29").Value = '0.997'
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("D30").Value = '1.71'
$ws.Range("E30").Value = '  +4.04%  '
$ws.Range("D31").Value = '6.11'
$ws.Range("E31").Value = '  +3.71%  '
$ws.Range("D32").Value = '157.75'
$ws.Range("E32").Value = '  +4.75%  '
$ws.Range("D33").Value = '19.14'
$ws.Range("E33").Value = '  +2.06%  '
$ws.Range("D34").Value = '4.08'
$ws.Range("E34").Value = '  +4.87%  '
$ws.Range("D35").Value = '0.915'
$ws.Range("E35").Value = '  +10.84%  '
$ws.Range("D36").Value = '0.910'
$ws.Range("E36").Value = '  +12.00%  '
$ws.Range("D37").Value = '1.16'
$ws.Range("E37").Value = '  +5.30%  '
$ws.Range("D38").Value = '37.48'
$ws.Range("E38").Value = '  +1.08%  '
$ws.Range("E39").Value = '  +5.72%  '
$ws.Range("D40").Value = '303.80'
$ws.Range("E40").Value = '  +8.05%  '
$ws.Range("D41").Value = '3.66'
$ws.Range("E41").Value = '  +2.29%  '
$ws.Range("D42").Value = '0.995'
$ws.Range("E42").Value = '  -0.40%  '
$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").Value = '0.0978'
$ws.Range("E43").Value = '  +4.43%  '
$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").Value = '0.604'
$ws.Range("E44").Value = '  +0.81%  '
$ws.Range("D45").Value = '0.0547'
$ws.Range("E45").Value = '  +3.03%  '
$ws.Range("D46").Value = '19.30'
$ws.Range("E46").Value = '  +1.26%  '
$ws.Range("D47").Value = '10.63'
$ws.Range("E47").Value = '  -0.54%  '
$ws.Range("E48").Value = '  +5.02%  '
$ws.Range("D49").Value = '122.50'
$ws.Range("E49").Value = '  +10.21%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = '4.64'
$ws.Range("E50").Value = '  +4.09%  '
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '1.958.89'
$ws.Range("E51").Value = '  +1.06%  '

Write-Host "Applied cryptos update"
